# Remove the unused "SaveMaps" / "MapsInterval" / "DrawLoadedSp" map-related
# parameters from the Description sheet and from each of the three
# ParameterFile data sheets (LandType = 0, 2, 9).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Description sheet: delete the three rows whose "COLUMN" (A) value
#    is SaveMaps, MapsInterval or DrawLoadedSp.
# ---------------------------------------------------------------------
$namesToRemove = @("SaveMaps", "MapsInterval", "DrawLoadedSp")

$wsDesc = $wb.Worksheets.Item("Description")
$lastRow = $wsDesc.Cells.Item($wsDesc.Rows.Count, 1).End(-4162).Row  # xlUp

$rowsToDelete = New-Object System.Collections.ArrayList
for ($r = 1; $r -le $lastRow; $r++) {
    $label = $wsDesc.Cells.Item($r, 1).Value2
    if ($namesToRemove -contains $label) {
        [void]$rowsToDelete.Add($r)
    }
}
# Delete from the bottom up so earlier row numbers stay valid.
$sortedRows = $rowsToDelete | Sort-Object -Descending
foreach ($r in $sortedRows) {
    $wsDesc.Rows.Item($r).Delete()
}

# ---------------------------------------------------------------------
# 2) ParameterFile data sheets: delete the matching columns.
# ---------------------------------------------------------------------
$sheetNames = @(
    "ParameterFile LandType = 0",
    "ParameterFile LandType = 2",
    "ParameterFile LandType = 9"
)

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $lastCol = $ws.Cells.Item(1, $ws.Columns.Count).End(-4159).Column  # xlToLeft

    $colsToDelete = New-Object System.Collections.ArrayList
    for ($c = 1; $c -le $lastCol; $c++) {
        $label = $ws.Cells.Item(1, $c).Value2
        if ($namesToRemove -contains $label) {
            [void]$colsToDelete.Add($c)
        }
    }
    $sortedCols = $colsToDelete | Sort-Object -Descending
    foreach ($c in $sortedCols) {
        $ws.Columns.Item($c).Delete()
    }
}
